$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime (de-de): 2016-09-01 13:06:25 -> 2016-09-01 13:07:17
$wsOverview.Range("G2").Value = "2016-09-01 13:07:17"
$wsDeDe.Range("H2").Value     = "2016-09-01 13:07:17"

# Latest Handoff Datetime (zh-cn): 2016-09-01 13:06:21 -> 2016-09-01 13:07:10
$wsZhCn.Range("H2").Value = "2016-09-01 13:07:10"

# Column width changes caused by the shorter "Ready for handoff" status text.
# (Target raw OOXML width is 17.2159881591797; ColumnWidth below is the closest
# input that round-trips to that value through this host's column-width model.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.38
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.38
